$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in Aidan's clock in/out times (columns L and M) for rows 11-13
$ws.Range("L11").Value = 0.79166666666666663
$ws.Range("M11").Value = 0.83333333333333337

$ws.Range("L12").Value = 0.54166666666666663
$ws.Range("M12").Value = 0.875

$ws.Range("L13").Value = 0.79166666666666663
$ws.Range("M13").Value = 0.875

# Update selection to M13 (no longer at P9), and clear the frozen/top-left scroll to G1
$ws.Activate()
$ws.Range("M13").Select()
